$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.333.97"
$ws.Cells.Item(2, 5).Value = "  -3.52%  "

$ws.Cells.Item(3, 4).Value = "3.100.11"
$ws.Cells.Item(3, 5).Value = "  -1.88%  "

$ws.Cells.Item(4, 5).Value = "  +0.00%  "

$ws.Cells.Item(5, 4).Value = "'552.79"
$ws.Cells.Item(5, 5).Value = "  -2.81%  "

$ws.Cells.Item(6, 4).Value = "'138.36"
$ws.Cells.Item(6, 5).Value = "  -8.70%  "

$ws.Cells.Item(7, 5).Value = "  +0.01%  "

$ws.Cells.Item(8, 4).Value = "3.095.46"
$ws.Cells.Item(8, 5).Value = "  -1.73%  "

$ws.Cells.Item(9, 4).Value = "'0.498"
$ws.Cells.Item(9, 5).Value = "  -1.46%  "

$ws.Cells.Item(10, 5).Value = "  -0.43%  "

$ws.Cells.Item(11, 4).Value = "'6.57"
$ws.Cells.Item(11, 5).Value = "  -3.15%  "

$ws.Cells.Item(12, 4).Value = "'0.461"
$ws.Cells.Item(12, 5).Value = "  -1.91%  "

$ws.Cells.Item(13, 4).Value = "'35.21"
$ws.Cells.Item(13, 5).Value = "  -6.60%  "

$ws.Cells.Item(14, 5).Value = "  -3.54%  "

$ws.Cells.Item(15, 4).Value = "3.600.63"
$ws.Cells.Item(15, 5).Value = "  -1.87%  "

$ws.Cells.Item(16, 4).Value = "63.303.26"
$ws.Cells.Item(16, 5).Value = "  -3.30%  "

$ws.Cells.Item(17, 5).Value = "  -1.04%  "

$ws.Cells.Item(18, 4).Value = "3.099.65"
$ws.Cells.Item(18, 5).Value = "  -1.85%  "

$ws.Cells.Item(19, 4).Value = "'506.91"
$ws.Cells.Item(19, 5).Value = "  -3.62%  "

$ws.Cells.Item(20, 4).Value = "'6.71"
$ws.Cells.Item(20, 5).Value = "  -2.66%  "

$ws.Cells.Item(21, 4).Value = "'13.60"
$ws.Cells.Item(21, 5).Value = "  -3.92%  "

$ws.Cells.Item(22, 4).Value = "'0.708"
$ws.Cells.Item(22, 5).Value = "  -0.87%  "

$ws.Cells.Item(23, 4).Value = "'7.27"
$ws.Cells.Item(23, 5).Value = "  -3.44%  "

$ws.Cells.Item(24, 4).Value = "'78.08"
$ws.Cells.Item(24, 5).Value = "  -2.58%  "

$ws.Cells.Item(25, 4).Value = "'12.37"
$ws.Cells.Item(25, 5).Value = "  -4.67%  "

$ws.Cells.Item(26, 5).Value = "  +0.05%  "

$ws.Cells.Item(27, 4).Value = "'2.76"
$ws.Cells.Item(27, 5).Value = "  -2.84%  "

$ws.Cells.Item(28, 4).Value = "'8.28"
$ws.Cells.Item(28, 5).Value = "  -7.57%  "

$ws.Cells.Item(29, 4).Value = "'1.00"
$ws.Cells.Item(29, 5).Value = "  +0.05%  "

$ws.Cells.Item(30, 5).Value = "  -10.23%  "

$ws.Cells.Item(31, 4).Value = "'26.49"
$ws.Cells.Item(31, 5).Value = "  -2.09%  "

$ws.Cells.Item(32, 5).Value = "  -7.41%  "

$ws.Cells.Item(33, 5).Value = "  -2.81%  "

$ws.Cells.Item(34, 2).Value = "Bittensor"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(34, 4).Value = "'526.30"
$ws.Cells.Item(34, 5).Value = "  -9.83%  "

$ws.Cells.Item(35, 2).Value = "OKB"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(35, 4).Value = "'57.57"
$ws.Cells.Item(35, 5).Value = "  +7.91%  "

$ws.Cells.Item(36, 4).Value = "'6.01"
$ws.Cells.Item(36, 5).Value = "  -3.12%  "

$ws.Cells.Item(37, 4).Value = "'5.23"
$ws.Cells.Item(37, 5).Value = "  -8.20%  "

$ws.Cells.Item(38, 4).Value = "'0.0415"
$ws.Cells.Item(38, 5).Value = "  -3.23%  "

$ws.Cells.Item(39, 4).Value = "3.082.38"
$ws.Cells.Item(39, 5).Value = "  +0.60%  "

$ws.Cells.Item(40, 4).Value = "'0.0795"
$ws.Cells.Item(40, 5).Value = "  -4.64%  "

$ws.Cells.Item(42, 4).Value = "'2.76"
$ws.Cells.Item(42, 5).Value = "  -9.79%  "

$ws.Cells.Item(43, 4).Value = "'8.14"
$ws.Cells.Item(43, 5).Value = "  -3.06%  "

$ws.Cells.Item(44, 4).Value = "'2.77"
$ws.Cells.Item(44, 5).Value = "  +74.92%  "

$ws.Cells.Item(45, 4).Value = "'0.254"
$ws.Cells.Item(45, 5).Value = "  -2.65%  "

$ws.Cells.Item(46, 5).Value = "  +0.09%  "

$ws.Cells.Item(47, 2).Value = "Fetch.AI"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(47, 4).Value = "'2.06"
$ws.Cells.Item(47, 5).Value = "  -7.69%  "

$ws.Cells.Item(48, 2).Value = "Monero"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(48, 4).Value = "'122.97"
$ws.Cells.Item(48, 5).Value = "  +0.93%  "

$ws.Cells.Item(49, 4).Value = "'24.43"
$ws.Cells.Item(49, 5).Value = "  -6.88%  "

$ws.Cells.Item(50, 4).Value = "'0.108"
$ws.Cells.Item(50, 5).Value = "  -2.84%  "

$ws.Cells.Item(51, 4).Value = "0.0₃0510"
$ws.Cells.Item(51, 5).Value = "  -7.38%  "

Write-Host "Applied cryptos list update"
